$wb = $excel.ActiveWorkbook

# --- Rename the "Device" sheet to "Apparatus" ---
$ws = $wb.Worksheets.Item("Device")

# --- Update the "Device" -> "Apparatus" text that lives in the sheet's own
#     cells (these back shared-string entries). Touch them in the same
#     order the strings are expected to re-appear in the shared string
#     table: "Apparatus type" (used twice), "Apparatus parameters",
#     "Apparatus type with default values:", then the summary sentence.
$ws.Range("B9").Value2 = "Apparatus type"
$ws.Range("B56").Value2 = "Apparatus type"
$ws.Range("C56").Value2 = "Apparatus parameters"
$ws.Range("A8").Value2 = "Apparatus type with default values:"
$ws.Range("A1").Value2 = "This sheet summarizes the apparatuses connected to buses."

$ws.Name = "Apparatus"

# --- Update the view/selection on the "Advance" sheet (B6 -> B16) without
#     disturbing which sheet/tab is active. ---
$wsAdv = $wb.Worksheets.Item("Advance")
$wsAdv.Range("B16").Select()

# --- Re-select the Apparatus sheet's A1 cell and re-activate it so it
#     stays the active tab (matching the original tabSelected/activeTab
#     state) and so its stale "G27" selection is cleared. ---
$ws.Range("A1").Select()

Write-Host "Device -> Apparatus rename complete"
